# Adds a new weekly price record for "Vega Modelo de Temuco - Berenjena".
# A new row is inserted before the existing row 86, shifting all the
# subsequent rows down by one (old row 86 becomes row 87, ..., old row 191
# becomes row 192). The sheet's used range therefore grows from A1:R191 to
# A1:R192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 86; everything from the old row 86
# downward (including formatting) shifts down by one row.
$ws.Rows("86:86").Insert()

# Populate the newly inserted row 86 with the new weekly record.
$ws.Range("A86").Value = 10
$ws.Range("B86").Value = "Vega Modelo de Temuco"
$ws.Range("C86").Value = "La Araucanía"
$ws.Range("D86").Value = 44494
$ws.Range("E86").Value = 9
$ws.Range("F86").Value = 100112001
$ws.Range("G86").Value = "Berenjena"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 90
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = 10000
$ws.Range("N86").Value = "$/caja 60 unidades"
$ws.Range("O86").Value = "Región de Arica y Parinacota"
$ws.Range("P86").Value = 167
$ws.Range("Q86").Value = 60
$ws.Range("R86").Value = "Hortaliza"
